# Experiment 1 throughput/latency fix
# - Rewrites the measured latency values on Sheet1 (columns B:E, rows 2-11)
#   with the corrected figures.
# - Widens/raises the line chart's value-axis scale to fit the corrected
#   (higher) data range.
# - Restores the worksheet's last-saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected data (NewReno, Reno, Tahoe, Vegas) -------------------------

$newReno = @(
    65.378132476900007,
    65.377361702100004,
    65.372012043400005,
    65.404364512200004,
    65.496138097100001,
    65.550685905500004,
    65.892020270299994,
    72.517706666699993,
    71.169846681899998,
    79.074042553200002
)

$reno = @(
    65.378132476900007,
    65.377361702100004,
    65.372012043400005,
    65.404364512200004,
    65.496138097100001,
    65.566465765800004,
    65.981368373400002,
    71.078984474899997,
    71.701067103100002,
    79.159333333299998
)

$tahoe = @(
    65.378132476900007,
    65.377361702100004,
    65.372012043400005,
    65.404364512200004,
    65.496138097100001,
    65.514874312299995,
    65.921846743299994,
    72.7691390552,
    71.567791237099996,
    79.590141414100003
)

$vegas = @(
    65.1182454992,
    65.327476268400005,
    65.325792962400001,
    65.332595744700001,
    65.371566345800005,
    65.5843394929,
    66.0864955901,
    76.925689502200001,
    77.273437447500001,
    73.073777777800004
)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newReno[$i]
    $ws.Cells.Item($row, 3).Value = $reno[$i]
    $ws.Cells.Item($row, 4).Value = $tahoe[$i]
    $ws.Cells.Item($row, 5).Value = $vegas[$i]
}

# --- Chart value-axis rescale ---------------------------------------------

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.MinimumScale = 60
$valueAxis.MaximumScale = 80

# --- Restore saved selection -----------------------------------------------

[void]$ws.Range("O14").Select()
